$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$league = "España 2"
$season = "2023/2024"

$teams = @(
    @('AD Alcorcón', 42, 0.76, 1.15, 12.12, 14.33, 31, 0.85, 34.43, 46.6, 11.62, 43.52, 41, 68, 43, 1.27),
    @('Albacete', 42, 1.09, 1.02, 12.62, 12.38, 20, 0.86, 41.5, 39.07, 9.28, 45.39, 48, 73, 49, 1.57),
    @('Burgos', 42, 0.98, 0.91, 12.19, 12.29, 27, 0.84, 39.76, 45.02, 11.79, 44.16, 46, 76, 47, 1.62),
    @('CD Eldense', 42, 0.82, 0.88, 10.98, 12.17, 28, 0.85, 35.98, 40.26, 9.63, 48.18, 47, 76, 46, 1.07),
    @('Eibar', 44, 1.04, 0.96, 13.02, 11.82, 31, 0.84, 46.07, 37.73, 9.84, 46.41, 47, 73, 54, 1.62),
    @('Elche', 42, 1.04, 0.86, 15.52, 9.9, 34, 0.85, 46.5, 32.17, 8.41, 46.59, 48, 69, 58, 1.22),
    @('Espanyol', 46, 1.09, 0.83, 12.07, 10.96, 34, 0.87, 41.37, 34.74, 9.08, 47.23, 48, 70, 51, 1.58),
    @('FC Andorra', 42, 0.89, 1.02, 13.21, 10.45, 27, 0.89, 52.71, 30.95, 8.88, 46.11, 47, 69, 62, 1.3),
    @('FC Cartagena', 42, 0.82, 1.09, 10.17, 13.07, 27, 0.88, 34.67, 46.33, 12.99, 43.63, 44, 71, 44, 1.61),
    @('Huesca', 42, 0.76, 0.78, 9.6, 11.4, 27, 0.8, 31.62, 42.24, 12.01, 44.1, 45, 72, 45, 1.23),
    @('Leganés', 42, 0.92, 0.66, 10.14, 10.67, 22, 0.86, 36.02, 43.14, 12.04, 44.96, 43, 73, 46, 1.38),
    @('Levante UD', 42, 1.07, 1.03, 13.05, 13.12, 26, 0.87, 42.33, 40.12, 9.78, 44.8, 46, 70, 52, 1.79),
    @('Mirandés', 42, 0.89, 1.04, 10.19, 13.62, 26, 0.83, 32.36, 44.86, 12.21, 43.61, 42, 70, 45, 1.42),
    @('Racing Ferrol', 42, 0.87, 0.89, 11.83, 10.48, 36, 0.86, 40.24, 37.48, 9.43, 45.45, 45, 72, 52, 1.33),
    @('Racing Santander', 42, 1.17, 1.15, 12.36, 14.62, 17, 0.88, 34.57, 43.21, 8.76, 46.18, 45, 72, 47, 1.97),
    @('Real Oviedo', 46, 1.02, 0.81, 12.48, 11.07, 32, 0.87, 44.15, 35.85, 8.75, 45.37, 45, 69, 54, 1.47),
    @('Real Valladolid', 42, 0.99, 0.78, 13.86, 10.43, 37, 0.85, 43.69, 32.05, 8.73, 46.34, 49, 66, 55, 1.66),
    @('Real Zaragoza', 42, 0.89, 0.77, 11.14, 10.33, 27, 0.89, 39.67, 35.62, 8.98, 45.24, 47, 73, 51, 1.24),
    @('SD Amorebieta', 42, 0.79, 1.12, 9.98, 13.07, 33, 0.85, 36.48, 44.52, 12.31, 43.13, 42, 71, 46, 1.19),
    @('Sporting Gijón', 44, 0.89, 1.02, 11.7, 12.93, 25, 0.86, 40.39, 41.25, 10.47, 41.95, 41, 65, 54, 1.46),
    @('Tenerife', 42, 0.96, 0.83, 11, 11.21, 29, 0.87, 38.45, 38.45, 9.39, 45.78, 47, 76, 49, 1.2),
    @('Villarreal II', 42, 0.96, 1.14, 13.24, 12.38, 26, 0.85, 37.88, 40.95, 11.13, 43.65, 44, 67, 50, 1.22)
)

$startRow = 393

$r = $startRow
foreach ($team in $teams) {
    $ws.Cells.Item($r, 1).Value = $team[0]
    $r = $r + 1
}

$r = $startRow
foreach ($team in $teams) {
    $ws.Cells.Item($r, 3).Value = $league
    $r = $r + 1
}

$r = $startRow
foreach ($team in $teams) {
    $ws.Cells.Item($r, 4).Value = $season
    $r = $r + 1
}

$r = $startRow
foreach ($team in $teams) {
    $ws.Cells.Item($r, 5).Value = $team[1]
    $ws.Cells.Item($r, 6).Value = $team[2]
    $ws.Cells.Item($r, 7).Value = $team[3]
    $ws.Cells.Item($r, 8).Value = $team[4]
    $ws.Cells.Item($r, 9).Value = $team[5]
    $ws.Cells.Item($r, 10).Value = $team[6]
    $ws.Cells.Item($r, 11).Value = $team[7]
    $ws.Cells.Item($r, 12).Value = $team[8]
    $ws.Cells.Item($r, 13).Value = $team[9]
    $ws.Cells.Item($r, 14).Value = $team[10]
    $ws.Cells.Item($r, 15).Value = $team[11]
    $ws.Cells.Item($r, 16).Value = $team[12]
    $ws.Cells.Item($r, 17).Value = $team[13]
    $ws.Cells.Item($r, 18).Value = $team[14]
    $ws.Cells.Item($r, 19).Value = $team[15]
    $r = $r + 1
}

$ws.Range("D393:D414").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 387
$win.ScrollColumn = 1
